$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 1.1267104356470823
$ws.Range("E5").Value = 1.394840708087148
$ws.Range("H5").Value = 0.027386657566182946
$ws.Range("K5").Value = 5.099255354139058
$ws.Range("N5").Value = 0.7923208540286591
$ws.Range("Q5").Value = 4.51010424958433
$ws.Range("T5").Value = 1.259375965331009
$ws.Range("W5").Value = 11.366613164905965
$ws.Range("Z5").Value = 3.273445595678131
$ws.Range("AC5").Value = 1.3643718163127525
$ws.Range("AF5").Value = 21.187887220686118
$ws.Range("AI5").Value = 4.048265039185108
$ws.Range("AL5").Value = 5.103645454449996
$ws.Range("AO5").Value = 0.793170618412877
$ws.Range("AR5").Value = 11.297946352189092

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 5.8428957695909505
$ws.Range("E5").Value = 0.9784217216366143
$ws.Range("H5").Value = 0.12738078309691164
$ws.Range("K5").Value = 5.140334897233011
$ws.Range("N5").Value = 2.3204554938983177
$ws.Range("Q5").Value = 24.258047309456327
$ws.Range("T5").Value = 5.825894925999356
$ws.Range("W5").Value = 64.47340462188538
$ws.Range("Z5").Value = 23.388939757862477
$ws.Range("AC5").Value = 0.9735869705382898
$ws.Range("AF5").Value = 15.688872267264822
$ws.Range("AI5").Value = 6.24344812591858
$ws.Range("AL5").Value = 5.138687747344227
$ws.Range("AO5").Value = 2.328526726689646
$ws.Range("AR5").Value = 18.332485206585567

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 8.954585102573107
$ws.Range("E5").Value = 1.1241725031133079
$ws.Range("H5").Value = 0.12653368409950969
$ws.Range("K5").Value = 6.314277100882539
$ws.Range("N5").Value = 2.639481971166968
$ws.Range("Q5").Value = 39.52494474970688
$ws.Range("T5").Value = 9.333919825213524
$ws.Range("W5").Value = 119.48249304444079
$ws.Range("Z5").Value = 52.24269293567045
$ws.Range("AC5").Value = 1.1679226396164306
$ws.Range("AF5").Value = 20.72043499101095
$ws.Range("AI5").Value = 7.958191974252465
$ws.Range("AL5").Value = 6.323143486707767
$ws.Range("AO5").Value = 2.6441911330793615
$ws.Range("AR5").Value = 22.769093805721543

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 0.8070572946288957
$ws.Range("E5").Value = 1.2828561683654582
$ws.Range("H5").Value = 0.01846352325138602
$ws.Range("K5").Value = 3.3119060112729723
$ws.Range("N5").Value = 0.4689949224555239
$ws.Range("Q5").Value = 3.5252705165332237
$ws.Range("T5").Value = 0.7356310083311977
$ws.Range("W5").Value = 7.476947342609292
$ws.Range("Z5").Value = 1.2779591176921283
$ws.Range("AC5").Value = 1.207088467857929
$ws.Range("AF5").Value = 13.671330653012317
$ws.Range("AI5").Value = 2.6294644824256594
$ws.Range("AL5").Value = 3.3179169274251987
$ws.Range("AO5").Value = 0.4693334153369293
$ws.Range("AR5").Value = 6.379984643407525

# Sheet 5
$ws = $wb.Worksheets.Item(5)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 0.6797028189686438
$ws.Range("E5").Value = 1.2326123313032282
$ws.Range("H5").Value = 0.018744248498849574
$ws.Range("K5").Value = 2.7082821575397715
$ws.Range("N5").Value = 0.6791408519377392
$ws.Range("Q5").Value = 1.63137475697944
$ws.Range("T5").Value = 0.5269401994233736
$ws.Range("W5").Value = 3.8071183804476436
$ws.Range("Z5").Value = 1.1573550734131444
$ws.Range("AC5").Value = 0.9734491328158817
$ws.Range("AF5").Value = 7.0548865432258525
$ws.Range("AI5").Value = 2.0180298964630206
$ws.Range("AL5").Value = 2.7098886586464896
$ws.Range("AO5").Value = 0.6794617526942971
$ws.Range("AR5").Value = 5.086412815715016

# Sheet 6
$ws = $wb.Worksheets.Item(6)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 4.769300036515664
$ws.Range("E5").Value = 0.9919847365076179
$ws.Range("H5").Value = 0.08126858592887584
$ws.Range("K5").Value = 7.604822517264164
$ws.Range("N5").Value = 2.251310638452998
$ws.Range("Q5").Value = 19.162054154673754
$ws.Range("T5").Value = 4.9276800815263835
$ws.Range("W5").Value = 83.20680746805705
$ws.Range("Z5").Value = 28.135666193939116
$ws.Range("AC5").Value = 1.2250134072358954
$ws.Range("AF5").Value = 30.34793971992549
$ws.Range("AI5").Value = 7.2937362146335305
$ws.Range("AL5").Value = 7.674026891129357
$ws.Range("AO5").Value = 2.2688361437221345
$ws.Range("AR5").Value = 29.434704785837738

# Sheet 7
$ws = $wb.Worksheets.Item(7)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 0.23889540291373781
$ws.Range("E5").Value = 0.878388556178404
$ws.Range("H5").Value = 0.0047328906813413315
$ws.Range("K5").Value = 1.8173453119931011
$ws.Range("N5").Value = 0.06127419014153981
$ws.Range("Q5").Value = 0.9655240461889877
$ws.Range("T5").Value = 0.204887964145779
$ws.Range("W5").Value = 2.0079350431463068
$ws.Range("Z5").Value = 0.22711162055973397
$ws.Range("AC5").Value = 1.359643317916559
$ws.Range("AF5").Value = 5.9138308357400815
$ws.Range("AI5").Value = 1.38059875176416
$ws.Range("AL5").Value = 1.8176066204849886
$ws.Range("AO5").Value = 0.06123472846283407
$ws.Range("AR5").Value = 1.8201983760944864

# Sheet 8
$ws = $wb.Worksheets.Item(8)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 1.591012410325903
$ws.Range("E5").Value = 1.3199637673584717
$ws.Range("H5").Value = 0.01981601624778787
$ws.Range("K5").Value = 5.220052816803952
$ws.Range("N5").Value = 0.8555350328466328
$ws.Range("Q5").Value = 8.659047857569837
$ws.Range("T5").Value = 1.5309581150521214
$ws.Range("W5").Value = 17.50515034780391
$ws.Range("Z5").Value = 6.420032761072121
$ws.Range("AC5").Value = 1.3977051799490414
$ws.Range("AF5").Value = 22.847030491812756
$ws.Range("AI5").Value = 4.731751000371239
$ws.Range("AL5").Value = 5.226515976024809
$ws.Range("AO5").Value = 0.8573342139243323
$ws.Range("AR5").Value = 13.51042643778789

# Sheet 9
$ws = $wb.Worksheets.Item(9)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 0.5652289506247141
$ws.Range("E5").Value = 1.3459484414905558
$ws.Range("H5").Value = 0.01211536209952054
$ws.Range("K5").Value = 3.2465982758491005
$ws.Range("N5").Value = 0.420245128878707
$ws.Range("Q5").Value = 1.8963356374223164
$ws.Range("T5").Value = 0.43159571988265594
$ws.Range("W5").Value = 4.309883921967491
$ws.Range("Z5").Value = 0.8072481719421457
$ws.Range("AC5").Value = 1.0863556158094254
$ws.Range("AF5").Value = 10.188057703132296
$ws.Range("AI5").Value = 1.6040353213140321
$ws.Range("AL5").Value = 3.2479779093504177
$ws.Range("AO5").Value = 0.42031037852831155
$ws.Range("AR5").Value = 4.270132497962854

# Sheet 10
$ws = $wb.Worksheets.Item(10)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 0.24721310380245634
$ws.Range("E5").Value = 1.116119567459998
$ws.Range("H5").Value = 0.005259563781964222
$ws.Range("K5").Value = 2.671025916319976
$ws.Range("N5").Value = 0.040957435631084456
$ws.Range("Q5").Value = 1.411935561573671
$ws.Range("T5").Value = 0.22186560668024377
$ws.Range("W5").Value = 3.3470516820028067
$ws.Range("Z5").Value = 0.25329333707696916
$ws.Range("AC5").Value = 1.2364106560868329
$ws.Range("AF5").Value = 11.170169837937863
$ws.Range("AI5").Value = 1.4356784137464385
$ws.Range("AL5").Value = 2.6724274805177646
$ws.Range("AO5").Value = 0.04094516673840532
$ws.Range("AR5").Value = 2.8919217984889407

# Sheet 11
$ws = $wb.Worksheets.Item(11)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 0.8382280030529543
$ws.Range("E5").Value = 1.1057723833851696
$ws.Range("H5").Value = 0.017714313710554788
$ws.Range("K5").Value = 2.6127111352173658
$ws.Range("N5").Value = 0.6035957918432372
$ws.Range("Q5").Value = 2.101363135251714
$ws.Range("T5").Value = 0.6526496603509886
$ws.Range("W5").Value = 4.48572893811848
$ws.Range("Z5").Value = 1.4022729412347548
$ws.Range("AC5").Value = 1.1125255342224292
$ws.Range("AF5").Value = 9.090860883127664
$ws.Range("AI5").Value = 2.597538674056246
$ws.Range("AL5").Value = 2.6142838971182596
$ws.Range("AO5").Value = 0.6040111984063129
$ws.Range("AR5").Value = 5.885256823750043

# Sheet 12
$ws = $wb.Worksheets.Item(12)
$ws.Range("A5").Value = "Average"
$ws.Range("B5").Value = 0.3147319037828517
$ws.Range("E5").Value = 1.348833102459005
$ws.Range("H5").Value = 0.0068649877147410636
$ws.Range("K5").Value = 2.4192791889348424
$ws.Range("N5").Value = 0.18767821701167017
$ws.Range("Q5").Value = 1.7758426106158105
$ws.Range("T5").Value = 0.2586503642120012
$ws.Range("W5").Value = 3.1611374602805156
$ws.Range("Z5").Value = 0.4022031929044953
$ws.Range("AC5").Value = 1.3001925134262595
$ws.Range("AF5").Value = 8.791430177359915
$ws.Range("AI5").Value = 1.6306574741579707
$ws.Range("AL5").Value = 2.4202334990122707
$ws.Range("AO5").Value = 0.18770187000984653
$ws.Range("AR5").Value = 3.1903303989662177
